$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells I1 and J1, copying the style of the existing header cell H1
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0

# Fill in the numeric data for rows 2 and 3
$ws.Range("I2").Value = 9
$ws.Range("J2").Value = 9
$ws.Range("I3").Value = 9
$ws.Range("J3").Value = 9
